$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Resumen")
$ws2 = $wb.Worksheets.Item("Solucion")
$ws3 = $wb.Worksheets.Item("Metricas")

# --- Resumen sheet: best zone & its metric changed ---
$ws1.Range("B2").Value = "Z4"
$ws1.Range("C2").Value = 527.2300789274517

# --- Solucion sheet: Pedido -> Salida assignments randomized/reshuffled ---
$ws2.Range("B3").Value = "S021"
$ws2.Range("B4").Value = "S041"
$ws2.Range("B5").Value = "S061"
$ws2.Range("B6").Value = "S071"
$ws2.Range("B7").Value = "S051"
$ws2.Range("B8").Value = "S011"
$ws2.Range("A9").Value = "Pedido_45"
$ws2.Range("B9").Value = "S031"
$ws2.Range("A10").Value = "Pedido_12"
$ws2.Range("B10").Value = "S062"
$ws2.Range("B11").Value = "S002"
$ws2.Range("A12").Value = "Pedido_76"
$ws2.Range("B12").Value = "S042"
$ws2.Range("B13").Value = "S022"
$ws2.Range("A14").Value = "Pedido_26"
$ws2.Range("B14").Value = "S072"
$ws2.Range("B15").Value = "S052"
$ws2.Range("B16").Value = "S012"
$ws2.Range("A17").Value = "Pedido_64"
$ws2.Range("B17").Value = "S032"
$ws2.Range("A18").Value = "Pedido_20"
$ws2.Range("B18").Value = "S003"
$ws2.Range("A19").Value = "Pedido_37"
$ws2.Range("B19").Value = "S043"
$ws2.Range("A20").Value = "Pedido_33"
$ws2.Range("B20").Value = "S063"
$ws2.Range("A21").Value = "Pedido_47"
$ws2.Range("B21").Value = "S023"
$ws2.Range("A22").Value = "Pedido_71"
$ws2.Range("B22").Value = "S053"
$ws2.Range("A23").Value = "Pedido_21"
$ws2.Range("B23").Value = "S033"
$ws2.Range("A24").Value = "Pedido_15"
$ws2.Range("B24").Value = "S013"
$ws2.Range("B25").Value = "S073"
$ws2.Range("B26").Value = "S024"
$ws2.Range("A27").Value = "Pedido_78"
$ws2.Range("B27").Value = "S044"
$ws2.Range("A28").Value = "Pedido_13"
$ws2.Range("B28").Value = "S004"
$ws2.Range("B29").Value = "S064"
$ws2.Range("B30").Value = "S054"
$ws2.Range("B31").Value = "S034"
$ws2.Range("A32").Value = "Pedido_35"
$ws2.Range("B32").Value = "S014"
$ws2.Range("A33").Value = "Pedido_6"
$ws2.Range("B33").Value = "S074"
$ws2.Range("A34").Value = "Pedido_52"
$ws2.Range("B34").Value = "S045"
$ws2.Range("A35").Value = "Pedido_4"
$ws2.Range("B35").Value = "S025"
$ws2.Range("B36").Value = "S005"
$ws2.Range("A37").Value = "Pedido_80"
$ws2.Range("B37").Value = "S065"
$ws2.Range("A38").Value = "Pedido_25"
$ws2.Range("B38").Value = "S055"
$ws2.Range("B39").Value = "S035"
$ws2.Range("B40").Value = "S075"
$ws2.Range("A41").Value = "Pedido_77"
$ws2.Range("B41").Value = "S015"
$ws2.Range("A42").Value = "Pedido_8"
$ws2.Range("B42").Value = "S046"
$ws2.Range("B43").Value = "S006"
$ws2.Range("A44").Value = "Pedido_39"
$ws2.Range("B44").Value = "S026"
$ws2.Range("A45").Value = "Pedido_23"
$ws2.Range("B45").Value = "S066"
$ws2.Range("B46").Value = "S056"
$ws2.Range("B47").Value = "S016"
$ws2.Range("B48").Value = "S036"
$ws2.Range("B49").Value = "S076"
$ws2.Range("A50").Value = "Pedido_53"
$ws2.Range("B50").Value = "S047"
$ws2.Range("A51").Value = "Pedido_51"
$ws2.Range("B51").Value = "S027"
$ws2.Range("B52").Value = "S007"
$ws2.Range("B53").Value = "S067"
$ws2.Range("A54").Value = "Pedido_55"
$ws2.Range("B54").Value = "S057"
$ws2.Range("A55").Value = "Pedido_18"
$ws2.Range("B55").Value = "S037"
$ws2.Range("B56").Value = "S017"
$ws2.Range("B57").Value = "S077"
$ws2.Range("B58").Value = "S048"
$ws2.Range("B59").Value = "S008"
$ws2.Range("B60").Value = "S028"
$ws2.Range("B61").Value = "S068"
$ws2.Range("A62").Value = "Pedido_70"
$ws2.Range("B62").Value = "S018"
$ws2.Range("A63").Value = "Pedido_54"
$ws2.Range("B63").Value = "S058"
$ws2.Range("A64").Value = "Pedido_67"
$ws2.Range("B64").Value = "S038"
$ws2.Range("A65").Value = "Pedido_40"
$ws2.Range("B65").Value = "S078"
$ws2.Range("B66").Value = "S009"
$ws2.Range("B67").Value = "S049"
$ws2.Range("B68").Value = "S029"
$ws2.Range("B69").Value = "S069"
$ws2.Range("B70").Value = "S059"
$ws2.Range("B71").Value = "S019"
$ws2.Range("A72").Value = "Pedido_66"
$ws2.Range("B72").Value = "S039"
$ws2.Range("A73").Value = "Pedido_62"
$ws2.Range("B73").Value = "S079"
$ws2.Range("A74").Value = "Pedido_63"
$ws2.Range("B74").Value = "S010"
$ws2.Range("B75").Value = "S050"
$ws2.Range("A76").Value = "Pedido_61"
$ws2.Range("B76").Value = "S030"
$ws2.Range("A77").Value = "Pedido_38"
$ws2.Range("B77").Value = "S070"
$ws2.Range("B78").Value = "S020"
$ws2.Range("B79").Value = "S060"
$ws2.Range("B80").Value = "S040"

# --- Metricas sheet: per-zone average time values updated ---
$ws3.Range("B2").Value = 516.402346199589
$ws3.Range("B3").Value = 526.3008487404044
$ws3.Range("B4").Value = 516.5745972537571
$ws3.Range("B5").Value = 527.2300789274517
